$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 173, pushing existing rows 173..252 down to 174..253
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(173,1).Value  = 3
$ws.Cells.Item(173,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(173,3).Value  = "Coquimbo"
$ws.Cells.Item(173,4).Value  = 44523
$ws.Cells.Item(173,5).Value  = 5
$ws.Cells.Item(173,6).Value  = 100112031
$ws.Cells.Item(173,7).Value  = "Poroto verde"
$ws.Cells.Item(173,8).Value  = "Magnum"
$ws.Cells.Item(173,9).Value  = "Primera"
$ws.Cells.Item(173,10).Value = 73
$ws.Cells.Item(173,11).Value = 39000
$ws.Cells.Item(173,12).Value = 40000
$ws.Cells.Item(173,13).Value = 39479
$ws.Cells.Item(173,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(173,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(173,16).Value = 1579
$ws.Cells.Item(173,17).Value = 25
$ws.Cells.Item(173,18).Value = "Hortaliza"
